$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C (pushes old C..I to D..J)
$ws.Columns("C").Insert()

# New header for the inserted column
$ws.Range("C1").Value = "Priority Number of the Problem"

# Renumber / update problem descriptions in column A
$ws.Range("A2").Value = "1. How does a quality of  case study can measure? Is it from length of our report?"
$ws.Range("A3").Value = "2. How should we divide the workload within our members to score a good mark for all?"
$ws.Range("A5").Value = "4. Since case study hasn't a right or wrong answer can we mention so many facts?"
$ws.Range("A4").Value = "3. Expectations for the Viva from sir"
$ws.Range("A6").Value = "5. Is there any option in twitter to measure the reliability of a tweet?"
$ws.Range("A7").Value = "6. If there's nothing for 6th problem, can we suggest to categorize the tweets and appoint an admin for each category?"
$ws.Range("A8").Value = "7. The Analytical Approach or The Problem-Oriented Method suitable for our academic case study?"

# Column width adjustments to match the new layout
$ws.Columns("A").ColumnWidth = 103.66666666666667
$ws.Columns("B").ColumnWidth = 51.666666666666664
$ws.Columns("C").ColumnWidth = 51.666666666666664
$ws.Columns("F").ColumnWidth = 10.666666666666666
$ws.Columns("J").ColumnWidth = 11.5
